$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: only E1 and F1 actually change text ("eventid" -> "id",
# "peerreviewid" -> "uniqueid"). The other header cells keep their text,
# but removing "eventid"/"peerreviewid" from the shared-string table and
# appending "uniqueid"/"id" at the end naturally happens as a side effect
# of changing these two cell values plus Excel's string-table bookkeeping.
$ws.Range("F1").Value = "uniqueid"
$ws.Range("E1").Value = "id"

# Move the active selection from K3 to E2.
$ws.Range("E2").Select()

$ws.Activate()
